# "fix the merge error"
#
# A previous merge clobbered the "View" column (F) flags for several
# Property rows and mixed up the Public/Private/Save/View flags on a
# couple of others. Restore the correct values and re-point the column F
# header at its own (un-clobbered) label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# --- Header (F1 = "View") -------------------------------------------------
# The merge left F1 sharing a shared-string slot with a stale copy; retype
# it so it gets its own distinct string entry again.
$f1 = $ws.Cells.Item(1, 6)
$f1.Value = "View"
$f1.Characters(1, 1).Font.Size = 11

# --- Rows 68-75: the View (F) flag was dropped by the merge -> restore it
for ($r = 68; $r -le 75; $r++) {
    $ws.Cells.Item($r, 6).Value = $true
}

# --- Rows 76-77: Private/Save were incorrectly left TRUE and View was
#     incorrectly left FALSE -> correct the flags
foreach ($r in 76, 77) {
    $ws.Cells.Item($r, 4).Value = $false   # Private
    $ws.Cells.Item($r, 5).Value = $false   # Save
    $ws.Cells.Item($r, 6).Value = $true    # View
}

# --- Row 78: Public flag was incorrectly left TRUE -> correct it
$ws.Cells.Item(78, 3).Value = $false

# Restore the selection that was active when the fix was made
$ws.Range("C78").Select()
